$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename columns to snake_case field names ---
$ws.Range("A1").Value = 'mx_state'
$ws.Range("B1").Value = 'mx_municipality'
$ws.Range("C1").Value = 'n_matriculas'
$ws.Range("D1").Value = 'pct_matriculas'

# --- Title-case municipality/state names (columns A and B) ---
$ws.Range("B5").Value = 'Pabellón De Arteaga'
$ws.Range("B6").Value = 'Rincón De Romos'
$ws.Range("B10").Value = 'Playas De Rosarito'
$ws.Range("B20").Value = 'Amatenango Del Valle'
$ws.Range("B23").Value = 'Benemérito De Las Américas'
$ws.Range("B30").Value = 'Comitán De Domínguez'
$ws.Range("B39").Value = 'Marqués De Comillas'
$ws.Range("B41").Value = 'Montecristo De Guerrero'
$ws.Range("B45").Value = 'Ocozocoautla De Espinosa'
$ws.Range("B50").Value = 'Salto De Agua'
$ws.Range("B76").Value = 'Hidalgo Del Parral'
$ws.Range("B82").Value = 'San Francisco Del Oro'
$ws.Range("B94").Value = 'San Juan De Sabinas'
$ws.Range("B104").Value = 'Villa De Álvarez'
$ws.Range("A106").Value = 'Ciudad De México'
$ws.Range("B134").Value = 'San Juan Del Río'
$ws.Range("B135").Value = 'San Pedro Del Gallo'
$ws.Range("A139").Value = 'Estado De México'
$ws.Range("B139").Value = 'Acambay De Ruíz Castañeda'
$ws.Range("B142").Value = 'Almoloya De Alquisiras'
$ws.Range("B146").Value = 'Atizapán De Zaragoza'
$ws.Range("B153").Value = 'Ecatepec De Morelos'
$ws.Range("B159").Value = 'Naucalpan De Juárez'
$ws.Range("B162").Value = 'San Felipe Del Progreso'
$ws.Range("B163").Value = 'San Simón De Guerrero'
$ws.Range("B169").Value = 'Tenango Del Valle'
$ws.Range("B170").Value = 'Tlalnepantla De Baz'
$ws.Range("B174").Value = 'Valle De Bravo'
$ws.Range("B181").Value = 'Apaseo El Alto'
$ws.Range("B182").Value = 'Apaseo El Grande'
$ws.Range("B188").Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range("B200").Value = 'San Diego De La Unión'
$ws.Range("B202").Value = 'San Francisco Del Rincón'
$ws.Range("B204").Value = 'San Luis De La Paz'
$ws.Range("B205").Value = 'Santa Cruz De Juventino Rosas'
$ws.Range("B206").Value = 'Silao De La Victoria'
$ws.Range("B210").Value = 'Valle De Santiago'
$ws.Range("B216").Value = 'Acapulco De Juárez'
$ws.Range("B220").Value = 'Atoyac De Álvarez'
$ws.Range("B221").Value = 'Ayutla De Los Libres'
$ws.Range("B222").Value = 'Buenavista De Cuéllar'
$ws.Range("B223").Value = 'Chilapa De Álvarez'
$ws.Range("B224").Value = 'Chilpancingo De Los Bravo'
$ws.Range("B227").Value = 'Coyuca De Benítez'
$ws.Range("B228").Value = 'Coyuca De Catalán'
$ws.Range("B230").Value = 'Cuetzala Del Progreso'
$ws.Range("B231").Value = 'Cutzamala De Pinzón'
$ws.Range("B234").Value = 'Huitzuco De Los Figueroa'
$ws.Range("B235").Value = 'Iguala De La Independencia'
$ws.Range("B236").Value = 'Zihuatanejo De Azueta'
$ws.Range("B238").Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range("B249").Value = 'Taxco De Alarcón'
$ws.Range("B253").Value = 'Tixtla De Guerrero'
$ws.Range("B256").Value = 'Tlapa De Comonfort'
$ws.Range("B262").Value = 'Atotonilco El Grande'
$ws.Range("B267").Value = 'Huejutla De Reyes'
$ws.Range("B270").Value = 'Jacala De Ledezma'
$ws.Range("B274").Value = 'Mixquiahuala De Juárez'
$ws.Range("B276").Value = 'Pachuca De Soto'
$ws.Range("B280").Value = 'Santiago De Anaya'
$ws.Range("B284").Value = 'Tenango De Doria'
$ws.Range("B285").Value = 'Tepehuacán De Guerrero'
$ws.Range("B286").Value = 'Tepeji Del Río De Ocampo'
$ws.Range("B287").Value = 'Tezontepec De Aldama'
$ws.Range("B293").Value = 'Tula De Allende'
$ws.Range("B294").Value = 'Tulancingo De Bravo'
$ws.Range("B303").Value = 'Encarnación De Díaz'
$ws.Range("B306").Value = 'Jilotlán De Los Dolores'
$ws.Range("B308").Value = 'Lagos De Moreno'
$ws.Range("B310").Value = 'Ojuelos De Jalisco'
$ws.Range("B313").Value = 'Talpa De Allende'
$ws.Range("B354").Value = 'Tiquicheo De Nicolás Romero'
$ws.Range("B370").Value = 'Puente De Ixtla'
$ws.Range("B373").Value = 'Tlaltizapán De Zapata'
$ws.Range("B375").Value = 'Zacualpan De Amilpas'
$ws.Range("B378").Value = 'Santa María Del Oro'
$ws.Range("B388").Value = 'Mier Y Noriega'
$ws.Range("B389").Value = 'Montemorelos'
$ws.Range("B392").Value = 'San Nicolás De Los Garza'
$ws.Range("B397").Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range("B398").Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range("B400").Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range("B401").Value = 'Oaxaca De Juárez'
$ws.Range("B402").Value = 'Ocotlán De Morelos'
$ws.Range("B403").Value = 'Putla Villa De Guerrero'
$ws.Range("B424").Value = 'San Miguel El Grande'
$ws.Range("B447").Value = 'Santo Domingo De Morelos'
$ws.Range("B451").Value = 'Teotitlán De Flores Magón'
$ws.Range("B452").Value = 'Tlacolula De Matamoros'
$ws.Range("B453").Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range("B459").Value = 'Ayotoxco De Guerrero'
$ws.Range("B465").Value = 'Cuetzalan Del Progreso'
$ws.Range("B473").Value = 'Izúcar De Matamoros'
$ws.Range("B477").Value = 'Los Reyes De Juárez'
$ws.Range("B479").Value = 'Palmar De Bravo'
$ws.Range("B483").Value = 'San Salvador El Seco'
$ws.Range("B490").Value = 'Tepexi De Rodríguez'
$ws.Range("B491").Value = 'Teteles De Avila Castillo'
$ws.Range("B494").Value = 'Tlacotepec De Benito Juárez'
$ws.Range("B499").Value = 'Tuzamapan De Galeana'
$ws.Range("B507").Value = 'Amealco De Bonfil'
$ws.Range("B509").Value = 'Cadereyta De Montes'
$ws.Range("B512").Value = 'Jalpan De Serra'
$ws.Range("B513").Value = 'Landa De Matamoros'
$ws.Range("B515").Value = 'Pinal De Amoles'
$ws.Range("B517").Value = 'San Juan Del Río'
$ws.Range("B529").Value = 'Ciudad Del Maíz'
$ws.Range("B538").Value = 'Mexquitic De Carmona'
$ws.Range("B543").Value = 'San Ciro De Acosta'
$ws.Range("B552").Value = 'Tanquián De Escobedo'
$ws.Range("B556").Value = 'Villa De Arista'
$ws.Range("B557").Value = 'Villa De Arriaga'
$ws.Range("B558").Value = 'Villa De Guadalupe'
$ws.Range("B559").Value = 'Villa De Ramos'
$ws.Range("B560").Value = 'Villa De Reyes'
$ws.Range("B621").Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range("B634").Value = 'Amatlán De Los Reyes'
$ws.Range("B641").Value = 'Boca Del Río'
$ws.Range("B649").Value = 'Cosamaloapan De Carpio'
$ws.Range("B650").Value = 'Cosautlán De Carvajal'
$ws.Range("B659").Value = 'Ignacio De La Llave'
$ws.Range("B661").Value = 'Ixhuatlán De Madero'
$ws.Range("B662").Value = 'Ixhuatlán Del Sureste'
$ws.Range("B670").Value = 'Lerdo De Tejada'
$ws.Range("B671").Value = 'Martínez De La Torre'
$ws.Range("B678").Value = 'Ozuluama De Mascareñas'
$ws.Range("B682").Value = 'Poza Rica De Hidalgo'
$ws.Range("B687").Value = 'Sayula De Alemán'
$ws.Range("B730").Value = 'Nochistlán De Mejía'
$ws.Range("B731").Value = 'Noria De Ángeles'
$ws.Range("B738").Value = 'Tlaltenango De Sánchez Román'

# --- Correct a floating point rounding value ---
$ws.Range("D571").Value = 0.009686168151879115

# --- Remove trailing footer/metadata rows (746-750) ---
$ws.Range("A746:D750").EntireRow.Delete()
